$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 into the new I1/J1 header cells, then set their text
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I and J columns, rows 2-55
$iValues = @{
    2 = 5
    3 = 6
    4 = 5
    5 = 5
    6 = 8
    7 = 7
    8 = 9
    9 = 6
    10 = 8
    11 = 8
    12 = 7
    13 = 8
    14 = 7
    15 = 8
    16 = 7
    17 = 7
    18 = 7
    19 = 7
    20 = 7
    21 = 9
    22 = 6
    23 = 7
    24 = 7
    25 = 8
    26 = 8
    27 = 8
    28 = 8
    29 = 8
    30 = 7
    31 = 7
    32 = 8
    33 = 8
    34 = 8
    35 = 6
    36 = 8
    37 = 7
    38 = 8
    39 = 8
    40 = 10
    41 = 8
    42 = 8
    43 = 8
    44 = 9
    45 = 8
    46 = 8
    47 = 9
    48 = 8
    49 = 8
    50 = 7
    51 = 8
    52 = 8
    53 = 9
    54 = 8
    55 = 6
}
$jValues = @{
    2 = 6
    3 = 6
    4 = 5
    5 = 5
    6 = 8
    7 = 8
    8 = 9
    9 = 6
    10 = 8
    11 = 8
    12 = 7
    13 = 8
    14 = 7
    15 = 8
    16 = 7
    17 = 7
    18 = 7
    19 = 7
    20 = 7
    21 = 9
    22 = 6
    23 = 7
    24 = 7
    25 = 8
    26 = 8
    27 = 8
    28 = 8
    29 = 8
    30 = 7
    31 = 7
    32 = 8
    33 = 8
    34 = 8
    35 = 6
    36 = 8
    37 = 7
    38 = 8
    39 = 8
    40 = 10
    41 = 8
    42 = 8
    43 = 8
    44 = 9
    45 = 8
    46 = 8
    47 = 9
    48 = 8
    49 = 8
    50 = 7
    51 = 8
    52 = 8
    53 = 9
    54 = 8
    55 = 6
}

foreach ($r in $iValues.Keys) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
    $ws.Cells.Item($r, 10).Value = $jValues[$r]
}
